# Update "model_input_variables_iran_se_calibrated.xlsx"
#  1. Update row 6 of sheet "strategy_id-0" (K6:AS6 all become 1)
#  2. Update row 2 of sheet "strategy_id-6003" (K2:U2 -> 1, V2:AS2 -> new decay curve)
#  3. Update row 2 of sheet "strategy_id-6004" (K2:U2 -> 1, V2:AS2 -> new decay curve)
#  4. Update row 2 of sheet "strategy_id-6005" (K2:U2 -> 1, V2:AS2 -> new decay curve, same as #3)
#  5. Add a new sheet "strategy_id-7032" at the end that is a duplicate of the
#     (now updated) "strategy_id-6005" sheet

$wb = $excel.ActiveWorkbook

# Columns K..AS (11..45), 35 contiguous columns
$cols = 11..45

# --- 1. strategy_id-0 : row 6, K6:AS6 -> 1 --------------------------------
$ws1 = $wb.Worksheets.Item("strategy_id-0")
$rowVals1 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws1.Cells.Item(6, $cols[$i]).Value = $rowVals1[$i]
}

# --- 2. strategy_id-6003 : row 2, K2:AS2 ----------------------------------
$ws2 = $wb.Worksheets.Item("strategy_id-6003")
$rowVals2 = @(1,1,1,1,1,1,1,1,1,1,1,0.9958333333333333,0.9916666666666666,0.9875,0.9833333333333334,0.9791666666666666,0.975,0.9708333333333332,0.9666666666666668,0.9625,0.9583333333333333,0.9541666666666667,0.95,0.9458333333333333,0.9416666666666667,0.9375,0.9333333333333333,0.9291666666666667,0.925,0.9208333333333334,0.9166666666666666,0.9125,0.9083333333333333,0.9041666666666667,0.9)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws2.Cells.Item(2, $cols[$i]).Value = $rowVals2[$i]
}

# --- 3. strategy_id-6004 : row 2, K2:AS2 ----------------------------------
$ws3 = $wb.Worksheets.Item("strategy_id-6004")
$rowVals3 = @(1,1,1,1,1,1,1,1,1,1,1,0.9875,0.975,0.9625,0.9500000000000001,0.9375,0.925,0.9124999999999999,0.9,0.8875,0.875,0.8625,0.85,0.8374999999999999,0.825,0.8125,0.8,0.7875,0.7749999999999999,0.7625,0.75,0.7374999999999999,0.725,0.7124999999999999,0.7)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws3.Cells.Item(2, $cols[$i]).Value = $rowVals3[$i]
}

# --- 4. strategy_id-6005 : row 2, K2:AS2 (same new curve as #3) ----------
$ws4 = $wb.Worksheets.Item("strategy_id-6005")
$rowVals4 = @(1,1,1,1,1,1,1,1,1,1,1,0.9875,0.975,0.9625,0.9500000000000001,0.9375,0.925,0.9124999999999999,0.9,0.8875,0.875,0.8625,0.85,0.8374999999999999,0.825,0.8125,0.8,0.7875,0.7749999999999999,0.7625,0.75,0.7374999999999999,0.725,0.7124999999999999,0.7)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws4.Cells.Item(2, $cols[$i]).Value = $rowVals4[$i]
}

# --- 5. Add new sheet "strategy_id-7032" as a copy of the updated --------
#        strategy_id-6005 sheet, placed after it (i.e. at the end)
$ws4.Copy([System.Type]::Missing, $ws4)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "strategy_id-7032"

Write-Host "Workbook now has" $wb.Worksheets.Count "sheets"
